# Auto-generated edit script: applies country table updates described in commit
# "Update countries & provincias Spain" - reorders a few countries in the ranking
# table (Australia/Brasil, Mexico/Crucero.../Croacia, El Salvador/Islas Virgenes...)
# and refreshes several numeric stat cells plus the "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = "Datos actualizados a 28 de Marzo de 2020 a las 03:59"
$ws.Cells.Item(4, 5).Value = 99924
$ws.Cells.Item(4, 8).Value = 1696
$ws.Cells.Item(21, 1).Value = "Australia"
$ws.Cells.Item(21, 2).Value = 3573
$ws.Cells.Item(21, 3).Value = 195
$ws.Cells.Item(21, 4).Value = 170
$ws.Cells.Item(21, 5).Value = 3389
$ws.Cells.Item(21, 6).Value = 23
$ws.Cells.Item(21, 8).Value = 14
$ws.Cells.Item(22, 1).Value = "Brasil"
$ws.Cells.Item(22, 2).Value = 3477
$ws.Cells.Item(22, 3).Value = 60
$ws.Cells.Item(22, 4).Value = 6
$ws.Cells.Item(22, 5).Value = 3378
$ws.Cells.Item(22, 6).Value = 296
$ws.Cells.Item(22, 7).Value = 1
$ws.Cells.Item(22, 8).Value = 93
$ws.Cells.Item(48, 1).Value = "Mexico"
$ws.Cells.Item(48, 2).Value = 717
$ws.Cells.Item(48, 3).Value = 132
$ws.Cells.Item(48, 4).Value = 4
$ws.Cells.Item(48, 5).Value = 701
$ws.Cells.Item(48, 6).Value = 1
$ws.Cells.Item(48, 7).Value = 4
$ws.Cells.Item(48, 8).Value = 12
$ws.Cells.Item(49, 1).Value = "Crucero"
$ws.Cells.Item(49, 2).Value = 712
$ws.Cells.Item(49, 3).Value = 0
$ws.Cells.Item(49, 4).Value = 597
$ws.Cells.Item(49, 5).Value = 105
$ws.Cells.Item(49, 6).Value = 15
$ws.Cells.Item(49, 7).Value = 0
$ws.Cells.Item(49, 8).Value = 10
$ws.Cells.Item(50, 1).Value = "Argentina"
$ws.Cells.Item(50, 2).Value = 690
$ws.Cells.Item(50, 3).Value = 101
$ws.Cells.Item(50, 4).Value = 72
$ws.Cells.Item(50, 5).Value = 601
$ws.Cells.Item(50, 6).Value = 0
$ws.Cells.Item(50, 7).Value = 4
$ws.Cells.Item(50, 8).Value = 17
$ws.Cells.Item(51, 1).Value = "Peru"
$ws.Cells.Item(51, 2).Value = 635
$ws.Cells.Item(51, 4).Value = 16
$ws.Cells.Item(51, 5).Value = 608
$ws.Cells.Item(51, 6).Value = 21
$ws.Cells.Item(51, 8).Value = 11
$ws.Cells.Item(52, 1).Value = "Eslovenia"
$ws.Cells.Item(52, 2).Value = 632
$ws.Cells.Item(52, 4).Value = 10
$ws.Cells.Item(52, 5).Value = 613
$ws.Cells.Item(52, 8).Value = 9
$ws.Cells.Item(53, 1).Value = "Croacia"
$ws.Cells.Item(53, 2).Value = 586
$ws.Cells.Item(53, 4).Value = 37
$ws.Cells.Item(53, 5).Value = 546
$ws.Cells.Item(53, 6).Value = 14
$ws.Cells.Item(53, 8).Value = 3
$ws.Cells.Item(74, 6).Value = 8
$ws.Cells.Item(103, 4).Value = 13
$ws.Cells.Item(103, 5).Value = 86
$ws.Cells.Item(104, 4).Value = 3
$ws.Cells.Item(104, 5).Value = 91
$ws.Cells.Item(142, 1).Value = "El Salvador"
$ws.Cells.Item(142, 2).Value = 19
$ws.Cells.Item(142, 3).Value = 6
$ws.Cells.Item(142, 5).Value = 19
$ws.Cells.Item(143, 1).Value = "Islas Virgenes de los Estados Unidos"
$ws.Cells.Item(143, 4).Value = 0
$ws.Cells.Item(143, 5).Value = 17
$ws.Cells.Item(144, 1).Value = "Bermudas"
$ws.Cells.Item(144, 2).Value = 17
$ws.Cells.Item(144, 4).Value = 2
$ws.Cells.Item(144, 5).Value = 15
$ws.Cells.Item(145, 1).Value = "Etiopia"
$ws.Cells.Item(145, 4).Value = 0
$ws.Cells.Item(145, 5).Value = 16
$ws.Cells.Item(146, 1).Value = "Maldivas"
$ws.Cells.Item(146, 2).Value = 16
$ws.Cells.Item(146, 4).Value = 9
$ws.Cells.Item(146, 5).Value = 7
$ws.Cells.Item(147, 1).Value = "Nueva Caledonia"
$ws.Cells.Item(147, 2).Value = 15
$ws.Cells.Item(147, 5).Value = 15
